$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.323836333333334
$ws.Range("H2").Value = 9.971509000000001
$ws.Range("I2").Value = 0.1301205485574112
$ws.Range("J2").Value = 0.1301205485574111
$ws.Range("M2").Value = 3.267668333333333
$ws.Range("N2").Value = 9.803005000000001
$ws.Range("O2").Value = 0.9196059551936462
$ws.Range("P2").Value = 0.9196059551936462
$ws.Range("Q2").Value = 10.86119473161611
$ws.Range("R2").Value = 97.75075258454501
$ws.Range("S2").Value = 0.1196596313464593
$ws.Range("T2").Value = 0.1196596313464593
$ws.Range("G3").Value = 3.323836333333334
$ws.Range("H3").Value = 9.971509000000001
$ws.Range("I3").Value = 0.1301205485574112
$ws.Range("J3").Value = 0.1301205485574111
$ws.Range("O3").Value = 0.08039404480635376
$ws.Range("P3").Value = 0.08039404480635376
$ws.Range("Q3").Value = 0.9495103538343335
$ws.Range("R3").Value = 8.545593184509
$ws.Range("S3").Value = 0.01046091721095184
$ws.Range("T3").Value = 0.01046091721095184
$ws.Range("I4").Value = 0.05261798027253429
$ws.Range("J4").Value = 0.05261798027253427
$ws.Range("M4").Value = 3.267668333333333
$ws.Range("N4").Value = 9.803005000000001
$ws.Range("O4").Value = 0.9196059551936462
$ws.Range("P4").Value = 0.9196059551936462
$ws.Range("Q4").Value = 4.392035973258889
$ws.Range("R4").Value = 39.52832375933
$ws.Range("S4").Value = 0.04838780800888433
$ws.Range("T4").Value = 0.04838780800888431
$ws.Range("I5").Value = 0.05261798027253429
$ws.Range("J5").Value = 0.05261798027253427
$ws.Range("O5").Value = 0.08039404480635376
$ws.Range("P5").Value = 0.08039404480635376
$ws.Range("S5").Value = 0.004230172263649959
$ws.Range("T5").Value = 0.004230172263649959
$ws.Range("G6").Value = 2.591152
$ws.Range("H6").Value = 7.773455999999999
$ws.Range("I6").Value = 0.1014376418761593
$ws.Range("J6").Value = 0.1014376418761593
$ws.Range("M6").Value = 3.267668333333333
$ws.Range("N6").Value = 9.803005000000001
$ws.Range("O6").Value = 0.9196059551936462
$ws.Range("P6").Value = 0.9196059551936462
$ws.Range("Q6").Value = 8.467025337253332
$ws.Range("R6").Value = 76.20322803528001
$ws.Range("S6").Value = 0.09328265955011644
$ws.Range("T6").Value = 0.09328265955011644
$ws.Range("G7").Value = 2.591152
$ws.Range("H7").Value = 7.773455999999999
$ws.Range("I7").Value = 0.1014376418761593
$ws.Range("J7").Value = 0.1014376418761593
$ws.Range("O7").Value = 0.08039404480635376
$ws.Range("P7").Value = 0.08039404480635376
$ws.Range("Q7").Value = 0.7402066183839999
$ws.Range("R7").Value = 6.661859565456
$ws.Range("S7").Value = 0.008154982326042814
$ws.Range("T7").Value = 0.008154982326042814
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.128738
$ws.Range("H8").Value = 0.386214
$ws.Range("I8").Value = 0.005039796638658401
$ws.Range("J8").Value = 0.005039796638658401
$ws.Range("M8").Value = 3.267668333333333
$ws.Range("N8").Value = 9.803005000000001
$ws.Range("O8").Value = 0.9196059551936462
$ws.Range("P8").Value = 0.9196059551936462
$ws.Range("Q8").Value = 0.4206730858966666
$ws.Range("R8").Value = 3.78605777307
$ws.Range("S8").Value = 0.004634627001875185
$ws.Range("T8").Value = 0.004634627001875185
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.128738
$ws.Range("H9").Value = 0.386214
$ws.Range("I9").Value = 0.005039796638658401
$ws.Range("J9").Value = 0.005039796638658401
$ws.Range("O9").Value = 0.08039404480635376
$ws.Range("P9").Value = 0.08039404480635376
$ws.Range("Q9").Value = 0.036776198246
$ws.Range("R9").Value = 0.330985784214
$ws.Range("S9").Value = 0.0004051696367832145
$ws.Range("T9").Value = 0.0004051696367832145
$ws.Range("G10").Value = 16.542033
$ws.Range("H10").Value = 49.626099
$ws.Range("I10").Value = 0.6475825499073805
$ws.Range("J10").Value = 0.6475825499073804
$ws.Range("M10").Value = 3.267668333333333
$ws.Range("N10").Value = 9.803005000000001
$ws.Range("O10").Value = 0.9196059551936462
$ws.Range("P10").Value = 0.9196059551936462
$ws.Range("Q10").Value = 54.053877403055
$ws.Range("R10").Value = 486.4848966274951
$ws.Range("S10").Value = 0.5955207693743138
$ws.Range("T10").Value = 0.5955207693743136
$ws.Range("G11").Value = 16.542033
$ws.Range("H11").Value = 49.626099
$ws.Range("I11").Value = 0.6475825499073805
$ws.Range("J11").Value = 0.6475825499073804
$ws.Range("O11").Value = 0.08039404480635376
$ws.Range("P11").Value = 0.08039404480635376
$ws.Range("Q11").Value = 4.725512941011
$ws.Range("R11").Value = 42.529616469099
$ws.Range("S11").Value = 0.05206178053306677
$ws.Range("T11").Value = 0.05206178053306676
$ws.Range("G12").Value = 1.614436666666667
$ws.Range("H12").Value = 4.84331
$ws.Range("I12").Value = 0.06320148274785642
$ws.Range("J12").Value = 0.06320148274785641
$ws.Range("M12").Value = 3.267668333333333
$ws.Range("N12").Value = 9.803005000000001
$ws.Range("O12").Value = 0.9196059551936462
$ws.Range("P12").Value = 0.9196059551936462
$ws.Range("Q12").Value = 5.275443571838888
$ws.Range("R12").Value = 47.47899214655
$ws.Range("S12").Value = 0.05812045991199725
$ws.Range("T12").Value = 0.05812045991199724
$ws.Range("G13").Value = 1.614436666666667
$ws.Range("H13").Value = 4.84331
$ws.Range("I13").Value = 0.06320148274785642
$ws.Range("J13").Value = 0.06320148274785641
$ws.Range("O13").Value = 0.08039404480635376
$ws.Range("P13").Value = 0.08039404480635376
$ws.Range("Q13").Value = 0.4611912792566666
$ws.Range("R13").Value = 4.15072151331
$ws.Range("S13").Value = 0.005081022835859163
$ws.Range("T13").Value = 0.005081022835859162
